$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Find first empty row after the existing data (row 15 is the last used row -> new row 16)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 2).Value = "Melegnano (MI)"
$ws.Cells.Item($newRow, 3).Value = "Musicolepsia"
$ws.Cells.Item($newRow, 1).Value = "Plug & Play Jam Session (Apr 6, 2023)"
$ws.Cells.Item($newRow, 4).Value = 45.357970954483299
$ws.Cells.Item($newRow, 5).Value = 9.3146710101658599
$ws.Cells.Item($newRow, 6).Value = 2023
$ws.Cells.Item($newRow, 7).Value = "06/04/2023"
$ws.Cells.Item($newRow, 8).Value = '<iframe width="300" height="169" src="https://www.youtube.com/embed/playlist?list=PLhIw1_0YGPEStVIUkVyv2ZB4PlUeK02QW"></iframe>'

$ws.Cells.Item($newRow, 1).Select()

# The H12/H13 cells previously carried a stray fill-applying style that is no
# longer needed; restore them to the plain "Normal" style so the unused xf
# entry drops out of the style table, matching the edited workbook.
$ws.Cells.Item(12, 8).Style = "Normal"
$ws.Cells.Item(13, 8).Style = "Normal"
